# Fix "Excel file total marks error" on the quiz marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking" row): Right count 5 -> 4, Wrong penalty -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total" row): total right marks 130 -> 104, total wrong penalty -2 -> -4
$ws.Range("B12").Value = 104
$ws.Range("C12").Value = -4

# Update the displayed score summary text to match the corrected totals
$ws.Range("E12").Value = "100 / 112"
